# Auto-generated Excel COM-interop script to update cryptos.xlsx
# Applies price/volume refresh values and the Kaspa/PancakeSwap row swap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.905.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.406.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.406.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.479"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.20%  "
$ws.Range("E11").Value = "  +5.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.979.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.374.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.875.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.572"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.531.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.72%  "
$ws.Range("E29").Value = "  +8.38%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.52%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.58%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.158"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.57%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.401.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0796"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.76%  "
$ws.Range("E42").Value = "  +12.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +8.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.773"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.346.22"
$ws.Range("D51").Style = "Normal"
